# edit.ps1 -- applies the two changes described by the diff:
#   1. The table on Slide 5 gets a new built-in table style GUID.
#   2. The presentation's theme color scheme is swapped from the
#      "Red Violet"/"Integral" palette back to the default
#      "Office"/"Office Theme" palette (font scheme and format scheme
#      are already identical between the two themes, so only the
#      12 theme colors actually need to change).

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$s = $p.Slides.Item(5)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{8B89ED6C-21FC-4469-8E19-6F743D101E35}")

# --- 2. Theme colour scheme: Red Violet / Integral -> Office / Office Theme
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Colors(1).RGB  = RGB 0x00 0x00 0x00   # dk1      -> 000000
$cs.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      -> FFFFFF
$cs.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      -> 44546A
$cs.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  -> ED7D31
$cs.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  -> FFC000
$cs.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  -> 4472C4
$cs.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6  -> 70AD47
$cs.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hlink    -> 0563C1
$cs.Colors(12).RGB = RGB 0x95 0x4F 0x72   # folHlink -> 954F72
